# Update Work Week and Social Spending
# Refreshes the Seychelles GDP per Capita data (years 1950-2016) on the
# "Data" sheet with newly downloaded values, and extends the series from
# 2008 through 2016 (8 new rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New "Data" column values for years 1950..2016, in order.
$values = @("2632", "2778", "2821", "2868", "2992", "2979", "2949", "3008", "3027", "3102", "3257", "2995", "3174", "3384", "3424", "3351", "3351", "3277", "3470", "3379", "3537", "3987", "4149", "4468", "4457", "4506", "5149", "5491", "5906", "6709", "6424", "5928", "5786", "5638", "6030", "6585", "6601", "6827", "7114", "7755", "8239", "8570.67385288007", "9304.80806777371", "10009.9140748174", "10060.8248378732", "10135.0958230014", "10468.9905824266", "12233.0496953541", "13801.1872520655", "14235.7867662143", "14379.8320383995", "14221.1627064738", "14581.4113899369", "13835.9936483827", "13735.8975558394", "15172.5631757869", "16824.6125397852", "18832.8258637132", "18687.3909992838", "18743.9969153866", "20147.1781530807", "21545", "22130", "23036", "24250", "25690", "26624")

$firstYear = 1950
$firstRow = 2
$lastExistingRow = 60

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $firstRow + $i
    $year = $firstYear + $i

    if ($row -gt $lastExistingRow) {
        # Brand-new row: fill in the constant columns too.
        $ws.Cells.Item($row, 1).Value = 690
        $ws.Cells.Item($row, 2).Value = "Seychelles"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = $year
    }

    # Keep the Data column as text (matches the source workbook, which
    # stores these figures as shared strings rather than numbers).
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
